# Update automàtic: dades i banners [2026-02-07 15:48]
#
# Refreshes the meteocat daily-summary scrape: per-station
# DATA_EXTRACCIO timestamps move forward ~30 min, and the
# associated HUMITAT_MITJANA_DIA, PRESSIO_ATMOSFERICA,
# RADIACIO_GLOBAL, RATXA_VENT_MAX, TEMPERATURA_MAXIMA_DIA and
# TEMPERATURA_MITJANA_DIA readings are updated to the newly
# fetched values. All of these columns are stored as plain
# text (inline strings) in the source sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Dades_Meteo")

# Columns with a value ending in "%" get auto-coerced into a numeric
# percentage by the COM Value setter (mirrors real Excel behaviour).
# Force these specific cells to keep a Text number format so the
# literal "NN%" string is preserved, matching the source data.
$percentCells = @("H2", "H3", "H5", "H7", "H8", "H9", "H10", "H11", "H12", "H13", "H14", "H15", "H16", "H17", "H19", "H21", "H22", "H24", "H25", "H27", "H28", "H29", "H33", "H34", "H36")
foreach ($pc in $percentCells) {
  $ws.Range($pc).NumberFormat = "@"
}

$ws.Range('E2').Value = '2026-02-07 15:47:34'
$ws.Range('H2').Value = '82%'
$ws.Range('K2').Value = '8.4 MJ/m2'
$ws.Range('O2').Value = '-0.7 °C'
$ws.Range('E3').Value = '2026-02-07 15:47:36'
$ws.Range('H3').Value = '84%'
$ws.Range('K3').Value = '13.6 MJ/m2'
$ws.Range('L3').Value = '27.7 km/h - 73º 15:29 TU'
$ws.Range('O3').Value = '-4.9 °C'
$ws.Range('E4').Value = '2026-02-07 15:47:38'
$ws.Range('J4').Value = '1003.4 hPa'
$ws.Range('K4').Value = '9.9 MJ/m2'
$ws.Range('M4').Value = '15.4 °C 15:24 TU'
$ws.Range('O4').Value = '12.0 °C'
$ws.Range('E5').Value = '2026-02-07 15:47:40'
$ws.Range('H5').Value = '65%'
$ws.Range('J5').Value = '1003.4 hPa'
$ws.Range('K5').Value = '11.5 MJ/m2'
$ws.Range('O5').Value = '10.4 °C'
$ws.Range('E6').Value = '2026-02-07 15:47:43'
$ws.Range('J6').Value = '1005.0 hPa'
$ws.Range('K6').Value = '11.8 MJ/m2'
$ws.Range('O6').Value = '13.0 °C'
$ws.Range('E7').Value = '2026-02-07 15:47:45'
$ws.Range('H7').Value = '63%'
$ws.Range('J7').Value = '1004.5 hPa'
$ws.Range('K7').Value = '12.3 MJ/m2'
$ws.Range('O7').Value = '9.2 °C'
$ws.Range('E8').Value = '2026-02-07 15:47:47'
$ws.Range('H8').Value = '76%'
$ws.Range('K8').Value = '11.1 MJ/m2'
$ws.Range('M8').Value = '16.0 °C 15:21 TU'
$ws.Range('O8').Value = '8.6 °C'
$ws.Range('E9').Value = '2026-02-07 15:47:50'
$ws.Range('H9').Value = '87%'
$ws.Range('O9').Value = '3.5 °C'
$ws.Range('E10').Value = '2026-02-07 15:47:52'
$ws.Range('H10').Value = '84%'
$ws.Range('O10').Value = '10.3 °C'
$ws.Range('E11').Value = '2026-02-07 15:47:54'
$ws.Range('H11').Value = '85%'
$ws.Range('K11').Value = '8.2 MJ/m2'
$ws.Range('O11').Value = '3.1 °C'
$ws.Range('E12').Value = '2026-02-07 15:47:56'
$ws.Range('H12').Value = '55%'
$ws.Range('K12').Value = '12.0 MJ/m2'
$ws.Range('O12').Value = '12.3 °C'
$ws.Range('E13').Value = '2026-02-07 15:47:58'
$ws.Range('H13').Value = '67%'
$ws.Range('O13').Value = '11.3 °C'
$ws.Range('E14').Value = '2026-02-07 15:48:01'
$ws.Range('H14').Value = '62%'
$ws.Range('K14').Value = '8.3 MJ/m2'
$ws.Range('L14').Value = '34.2 km/h - 242º 15:22 TU'
$ws.Range('O14').Value = '-5.8 °C'
$ws.Range('E15').Value = '2026-02-07 15:48:03'
$ws.Range('H15').Value = '72%'
$ws.Range('J15').Value = '1003.7 hPa'
$ws.Range('K15').Value = '11.0 MJ/m2'
$ws.Range('O15').Value = '9.5 °C'
$ws.Range('E16').Value = '2026-02-07 15:48:05'
$ws.Range('H16').Value = '87%'
$ws.Range('K16').Value = '6.2 MJ/m2'
$ws.Range('O16').Value = '3.8 °C'
$ws.Range('E17').Value = '2026-02-07 15:48:08'
$ws.Range('H17').Value = '87%'
$ws.Range('K17').Value = '8.4 MJ/m2'
$ws.Range('M17').Value = '10.0 °C 15:09 TU'
$ws.Range('O17').Value = '4.8 °C'
$ws.Range('E18').Value = '2026-02-07 15:48:10'
$ws.Range('K18').Value = '5.8 MJ/m2'
$ws.Range('O18').Value = '-5.7 °C'
$ws.Range('E19').Value = '2026-02-07 15:48:13'
$ws.Range('H19').Value = '82%'
$ws.Range('K19').Value = '11.5 MJ/m2'
$ws.Range('O19').Value = '6.8 °C'
$ws.Range('E20').Value = '2026-02-07 15:48:15'
$ws.Range('K20').Value = '10.0 MJ/m2'
$ws.Range('O20').Value = '-3.7 °C'
$ws.Range('E21').Value = '2026-02-07 15:48:17'
$ws.Range('H21').Value = '68%'
$ws.Range('K21').Value = '11.0 MJ/m2'
$ws.Range('O21').Value = '8.3 °C'
$ws.Range('E22').Value = '2026-02-07 15:48:19'
$ws.Range('H22').Value = '74%'
$ws.Range('K22').Value = '12.3 MJ/m2'
$ws.Range('O22').Value = '10.3 °C'
$ws.Range('E23').Value = '2026-02-07 15:48:22'
$ws.Range('K23').Value = '9.8 MJ/m2'
$ws.Range('O23').Value = '10.4 °C'
$ws.Range('E24').Value = '2026-02-07 15:48:24'
$ws.Range('H24').Value = '72%'
$ws.Range('J24').Value = '1002.8 hPa'
$ws.Range('K24').Value = '9.5 MJ/m2'
$ws.Range('O24').Value = '11.1 °C'
$ws.Range('E25').Value = '2026-02-07 15:48:26'
$ws.Range('H25').Value = '86%'
$ws.Range('J25').Value = '1005.9 hPa'
$ws.Range('K25').Value = '6.7 MJ/m2'
$ws.Range('M25').Value = '7.1 °C 15:14 TU'
$ws.Range('O25').Value = '2.1 °C'
$ws.Range('E26').Value = '2026-02-07 15:48:29'
$ws.Range('K26').Value = '10.4 MJ/m2'
$ws.Range('E27').Value = '2026-02-07 15:48:31'
$ws.Range('H27').Value = '76%'
$ws.Range('K27').Value = '10.6 MJ/m2'
$ws.Range('O27').Value = '11.4 °C'
$ws.Range('E28').Value = '2026-02-07 15:48:34'
$ws.Range('H28').Value = '82%'
$ws.Range('M28').Value = '10.3 °C 15:00 TU'
$ws.Range('O28').Value = '4.5 °C'
$ws.Range('E29').Value = '2026-02-07 15:48:36'
$ws.Range('H29').Value = '53%'
$ws.Range('K29').Value = '11.9 MJ/m2'
$ws.Range('O29').Value = '12.3 °C'
$ws.Range('E30').Value = '2026-02-07 15:48:38'
$ws.Range('K30').Value = '13.6 MJ/m2'
$ws.Range('E31').Value = '2026-02-07 15:48:41'
$ws.Range('O31').Value = '5.4 °C'
$ws.Range('E32').Value = '2026-02-07 15:48:43'
$ws.Range('J32').Value = '1006.3 hPa'
$ws.Range('K32').Value = '11.6 MJ/m2'
$ws.Range('O32').Value = '13.0 °C'
$ws.Range('E33').Value = '2026-02-07 15:48:46'
$ws.Range('H33').Value = '77%'
$ws.Range('O33').Value = '10.2 °C'
$ws.Range('E34').Value = '2026-02-07 15:48:48'
$ws.Range('H34').Value = '72%'
$ws.Range('K34').Value = '9.6 MJ/m2'
$ws.Range('O34').Value = '7.6 °C'
$ws.Range('E35').Value = '2026-02-07 15:48:50'
$ws.Range('K35').Value = '5.7 MJ/m2'
$ws.Range('O35').Value = '-4.4 °C'
$ws.Range('E36').Value = '2026-02-07 15:48:53'
$ws.Range('H36').Value = '79%'
$ws.Range('K36').Value = '11.2 MJ/m2'
$ws.Range('O36').Value = '8.4 °C'
